$d = $word.ActiveDocument

# 1) Intro paragraph: "discretewq v2.3.1: https://github.com/..." -> v2.3.2
$r1 = $d.Content
$r1.Find.Execute("discretewq v2.3.1: https", $true, $false, $false, $false, $false,
                  $true, 1, $false, "discretewq v2.3.2: https", 2)

# 2) "R package v2.3.1 (https://github.com/..." -> v2.3.2
$r2 = $d.Content
$r2.Find.Execute("R package v2.3.1 (https", $true, $false, $false, $false, $false,
                  $true, 1, $false, "R package v2.3.2 (https", 2)

# 3) Zenodo citation: update version and DOI
$r3 = $d.Content
$r3.Find.Execute("Estuary v2.3.1. Zenodo. doi:10.5281/zenodo.6335814", $true, $false, $false, $false, $false,
                  $true, 1, $false, "Estuary v2.3.2. Zenodo. doi:10.5281/zenodo.6390964", 2)
